$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3.648214333333333
$ws.Range("H2").Value = 10.944643
$ws.Range("I2").Value = 0.1438182892240057
$ws.Range("J2").Value = 0.1438182892240057
$ws.Range("M2").Value = 2.229346333333333
$ws.Range("N2").Value = 6.688039
$ws.Range("O2").Value = 0.3059994343984824
$ws.Range("P2").Value = 0.3059994343984824
$ws.Range("Q2").Value = 8.133133247230777
$ws.Range("R2").Value = 73.19819922507699
$ws.Range("S2").Value = 0.04400831515870311
$ws.Range("T2").Value = 0.04400831515870309
$ws.Range("G3").Value = 3.648214333333333
$ws.Range("H3").Value = 10.944643
$ws.Range("I3").Value = 0.1438182892240057
$ws.Range("J3").Value = 0.1438182892240057
$ws.Range("O3").Value = 0.6940005656015176
$ws.Range("P3").Value = 0.6940005656015176
$ws.Range("Q3").Value = 18.44578270148155
$ws.Range("R3").Value = 166.012044313334
$ws.Range("S3").Value = 0.09980997406530261
$ws.Range("T3").Value = 0.09980997406530259
$ws.Range("H4").Value = 9.050426
$ws.Range("I4").Value = 0.1189272947567555
$ws.Range("J4").Value = 0.1189272947567555
$ws.Range("M4").Value = 2.229346333333333
$ws.Range("N4").Value = 6.688039
$ws.Range("O4").Value = 0.3059994343984824
$ws.Range("P4").Value = 0.3059994343984824
$ws.Range("Q4").Value = 6.725511339401556
$ws.Range("R4").Value = 60.529602054614
$ws.Range("S4").Value = 0.0363916849301088
$ws.Range("T4").Value = 0.03639168493010879
$ws.Range("H5").Value = 9.050426
$ws.Range("I5").Value = 0.1189272947567555
$ws.Range("J5").Value = 0.1189272947567555
$ws.Range("O5").Value = 0.6940005656015176
$ws.Range("P5").Value = 0.6940005656015176
$ws.Range("S5").Value = 0.08253560982664673
$ws.Range("T5").Value = 0.08253560982664673
$ws.Range("G6").Value = 10.36785533333333
$ws.Range("H6").Value = 31.103566
$ws.Range("I6").Value = 0.4087169998040092
$ws.Range("J6").Value = 0.4087169998040092
$ws.Range("M6").Value = 2.229346333333333
$ws.Range("N6").Value = 6.688039
$ws.Range("O6").Value = 0.3059994343984824
$ws.Range("P6").Value = 0.3059994343984824
$ws.Range("Q6").Value = 23.11354027189711
$ws.Range("R6").Value = 208.021862447074
$ws.Range("S6").Value = 0.1250671707690715
$ws.Range("T6").Value = 0.1250671707690715
$ws.Range("G7").Value = 10.36785533333333
$ws.Range("H7").Value = 31.103566
$ws.Range("I7").Value = 0.4087169998040092
$ws.Range("J7").Value = 0.4087169998040092
$ws.Range("O7").Value = 0.6940005656015176
$ws.Range("P7").Value = 0.6940005656015176
$ws.Range("Q7").Value = 52.42104467703422
$ws.Range("R7").Value = 471.789402093308
$ws.Range("S7").Value = 0.2836498290349377
$ws.Range("T7").Value = 0.2836498290349377
$ws.Range("G8").Value = 0.05357233333333333
$ws.Range("H8").Value = 0.160717
$ws.Range("I8").Value = 0.00211190479115806
$ws.Range("J8").Value = 0.00211190479115806
$ws.Range("M8").Value = 2.229346333333333
$ws.Range("N8").Value = 6.688039
$ws.Range("O8").Value = 0.3059994343984824
$ws.Range("P8").Value = 0.3059994343984824
$ws.Range("Q8").Value = 0.1194312848847778
$ws.Range("R8").Value = 1.074881563963
$ws.Range("S8").Value = 0.0006462416715978116
$ws.Range("T8").Value = 0.0006462416715978115
$ws.Range("G9").Value = 0.05357233333333333
$ws.Range("H9").Value = 0.160717
$ws.Range("I9").Value = 0.00211190479115806
$ws.Range("J9").Value = 0.00211190479115806
$ws.Range("O9").Value = 0.6940005656015176
$ws.Range("P9").Value = 0.6940005656015176
$ws.Range("Q9").Value = 0.2708677531495555
$ws.Range("R9").Value = 2.437809778346
$ws.Range("S9").Value = 0.001465663119560249
$ws.Range("T9").Value = 0.001465663119560249
$ws.Range("G10").Value = 8.280381
$ws.Range("H10").Value = 24.841143
$ws.Range("I10").Value = 0.3264255114240716
$ws.Range("J10").Value = 0.3264255114240716
$ws.Range("M10").Value = 2.229346333333333
$ws.Range("N10").Value = 6.688039
$ws.Range("O10").Value = 0.3059994343984824
$ws.Range("P10").Value = 0.3059994343984824
$ws.Range("Q10").Value = 18.459837020953
$ws.Range("R10").Value = 166.138533188577
$ws.Range("S10").Value = 0.09988602186900128
$ws.Range("T10").Value = 0.09988602186900124
$ws.Range("G11").Value = 8.280381
$ws.Range("H11").Value = 24.841143
$ws.Range("I11").Value = 0.3264255114240716
$ws.Range("J11").Value = 0.3264255114240716
$ws.Range("O11").Value = 0.6940005656015176
$ws.Range("P11").Value = 0.6940005656015176
$ws.Range("Q11").Value = 41.866539258926
$ws.Range("R11").Value = 376.798853330334
$ws.Range("S11").Value = 0.2265394895550703
$ws.Range("T11").Value = 0.2265394895550703
